$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 10.61535666666667
$ws.Range("H2").Value = 31.84607
$ws.Range("I2").Value = 0.1437966543677284
$ws.Range("J2").Value = 0.1437966543677284
$ws.Range("M2").Value = 4.230734666666667
$ws.Range("N2").Value = 12.692204
$ws.Range("O2").Value = 0.3081346507358854
$ws.Range("P2").Value = 0.3081346507358855
$ws.Range("Q2").Value = 44.91075744869778
$ws.Range("R2").Value = 404.19681703828
$ws.Range("S2").Value = 0.04430873187058883
$ws.Range("T2").Value = 0.04430873187058884
# Row 3
$ws.Range("G3").Value = 10.61535666666667
$ws.Range("H3").Value = 31.84607
$ws.Range("I3").Value = 0.1437966543677284
$ws.Range("J3").Value = 0.1437966543677284
$ws.Range("O3").Value = 0.6225996527787135
$ws.Range("P3").Value = 0.6225996527787135
$ws.Range("Q3").Value = 90.74416631434001
$ws.Range("R3").Value = 816.6974968290601
$ws.Range("S3").Value = 0.0895277470800884
$ws.Range("T3").Value = 0.0895277470800884
# Row 4
$ws.Range("G4").Value = 10.61535666666667
$ws.Range("H4").Value = 31.84607
$ws.Range("I4").Value = 0.1437966543677284
$ws.Range("J4").Value = 0.1437966543677284
$ws.Range("M4").Value = 0.9510283333333334
$ws.Range("N4").Value = 2.853085
$ws.Range("O4").Value = 0.0692656964854011
$ws.Range("P4").Value = 0.0692656964854011
$ws.Range("Q4").Value = 10.09550495843889
$ws.Range("R4").Value = 90.85954462595001
$ws.Range("S4").Value = 0.009960175417051203
$ws.Range("T4").Value = 0.009960175417051203
# Row 5
$ws.Range("I5").Value = 0.6785840820776819
$ws.Range("J5").Value = 0.6785840820776819
$ws.Range("M5").Value = 4.230734666666667
$ws.Range("N5").Value = 12.692204
$ws.Range("O5").Value = 0.3081346507358854
$ws.Range("P5").Value = 0.3081346507358855
$ws.Range("Q5").Value = 211.9362599410902
$ws.Range("R5").Value = 1907.426339469812
$ws.Range("S5").Value = 0.2090952691259379
$ws.Range("T5").Value = 0.209095269125938
# Row 6
$ws.Range("I6").Value = 0.6785840820776819
$ws.Range("J6").Value = 0.6785840820776819
$ws.Range("O6").Value = 0.6225996527787135
$ws.Range("P6").Value = 0.6225996527787135
$ws.Range("S6").Value = 0.4224862138827268
$ws.Range("T6").Value = 0.4224862138827268
# Row 7
$ws.Range("I7").Value = 0.6785840820776819
$ws.Range("J7").Value = 0.6785840820776819
$ws.Range("M7").Value = 0.9510283333333334
$ws.Range("N7").Value = 2.853085
$ws.Range("O7").Value = 0.0692656964854011
$ws.Range("P7").Value = 0.0692656964854011
$ws.Range("Q7").Value = 47.64122639330611
$ws.Range("R7").Value = 428.771037539755
$ws.Range("S7").Value = 0.04700259906901722
$ws.Range("T7").Value = 0.04700259906901722
# Row 8
$ws.Range("G8").Value = 13.11220933333333
$ws.Range("H8").Value = 39.336628
$ws.Range("I8").Value = 0.1776192635545896
$ws.Range("J8").Value = 0.1776192635545896
$ws.Range("M8").Value = 4.230734666666667
$ws.Range("N8").Value = 12.692204
$ws.Range("O8").Value = 0.3081346507358854
$ws.Range("P8").Value = 0.3081346507358855
$ws.Range("Q8").Value = 55.47427858312356
$ws.Range("R8").Value = 499.2685072481121
$ws.Range("S8").Value = 0.05473064973935865
$ws.Range("T8").Value = 0.05473064973935866
# Row 9
$ws.Range("G9").Value = 13.11220933333333
$ws.Range("H9").Value = 39.336628
$ws.Range("I9").Value = 0.1776192635545896
$ws.Range("J9").Value = 0.1776192635545896
$ws.Range("O9").Value = 0.6225996527787135
$ws.Range("P9").Value = 0.6225996527787135
$ws.Range("Q9").Value = 112.088226694136
$ws.Range("R9").Value = 1008.794040247224
$ws.Range("S9").Value = 0.1105856918158983
$ws.Range("T9").Value = 0.1105856918158983
# Row 10
$ws.Range("G10").Value = 13.11220933333333
$ws.Range("H10").Value = 39.336628
$ws.Range("I10").Value = 0.1776192635545896
$ws.Range("J10").Value = 0.1776192635545896
$ws.Range("M10").Value = 0.9510283333333334
$ws.Range("N10").Value = 2.853085
$ws.Range("O10").Value = 0.0692656964854011
$ws.Range("P10").Value = 0.0692656964854011
$ws.Range("Q10").Value = 12.47008258859778
$ws.Range("R10").Value = 112.23074329738
$ws.Range("S10").Value = 0.01230292199933267
$ws.Range("T10").Value = 0.01230292199933267
